# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.301.05"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.690.49"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'217.94"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'0.5353"
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'0.2725"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").Value = "'0.06427"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "'21.73"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").Value = "1.706.49"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "'4.528"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'0.5801"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "'66.88"
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("D17").Value = "26.317.15"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "'4.911"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "'193.44"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("D22").Value = "'6.277"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").Value = "'1.007"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'148.72"
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("D25").Value = "'0.1288"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").Value = "'7.871"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").Value = "'15.84"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").Value = "'1.383"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").Value = "'0.06107"
$ws.Range("E29").Value = "  -5.70%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'3.585"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "'1.687"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").Value = "'1.034"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").Value = "'0.6189"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "'2.425"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").Value = "'6.221"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").Value = "'0.01642"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("D40").Value = "1.112.13"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "'0.8773"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "'100.92"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "1.841.88"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").Value = "'57.85"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").Value = "'1.012"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").Value = "'8.146"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "'0.05291"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").Value = "'0.4291"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'6.059"
$ws.Range("E51").Value = "  -0.23%  "
